$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.047.99'
$ws.Range("E2").Value = '  +3.44%  '

$ws.Range("D3").Value = '2.361.77'
$ws.Range("E3").Value = '  +1.39%  '

$ws.Range("E4").Value = '  +0.40%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '309.69'
$ws.Range("E5").Value = '  -0.64%  '

$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '108.16'
$ws.Range("E6").Value = '  -0.94%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  +0.23%  '

$ws.Range("D9").Value = '0.615'
$ws.Range("E9").Value = '  -0.26%  '

$ws.Range("D10").Value = '40.89'
$ws.Range("E10").Value = '  -0.40%  '

$ws.Range("D11").Value = '0.0916'
$ws.Range("E11").Value = '  -0.43%  '

$ws.Range("D12").Value = '8.44'
$ws.Range("E12").Value = '  -1.94%  '

$ws.Range("E13").Value = '  +1.43%  '

$ws.Range("D14").Value = '0.975'
$ws.Range("E14").Value = '  -2.75%  '

$ws.Range("D15").Value = '2.725.23'
$ws.Range("E15").Value = '  +1.70%  '

$ws.Range("D16").Value = '15.21'
$ws.Range("E16").Value = '  -1.77%  '

$ws.Range("D17").Value = '2.368.23'
$ws.Range("E17").Value = '  +1.84%  '

$ws.Range("D18").Value = '45.092.13'
$ws.Range("E18").Value = '  +4.36%  '

$ws.Range("D19").Value = '14.73'
$ws.Range("E19").Value = '  +12.50%  '

$ws.Range("D20").Value = '7.26'
$ws.Range("E20").Value = '  -4.13%  '

$ws.Range("E21").Value = '  -1.03%  '

$ws.Range("D22").Value = '73.11'
$ws.Range("E22").Value = '  -1.40%  '

$ws.Range("D23").Value = '3.47'
$ws.Range("E23").Value = '  -0.88%  '

$ws.Range("D24").Value = '258.63'
$ws.Range("E24").Value = '  -3.69%  '

$ws.Range("D25").Value = '2.29'
$ws.Range("E25").Value = '  +0.70%  '

$ws.Range("E26").Value = '  -0.42%  '

$ws.Range("D27").Value = '11.09'
$ws.Range("E27").Value = '  -0.66%  '

$ws.Range("D28").Value = '7.22'
$ws.Range("E28").Value = '  -5.11%  '

$ws.Range("E29").Value = '  +2.05%  '

$ws.Range("D30").Value = '0.0963'
$ws.Range("E30").Value = '  +8.40%  '

$ws.Range("D31").Value = '22.35'
$ws.Range("E31").Value = '  -1.33%  '

$ws.Range("D32").Value = '37.32'
$ws.Range("E32").Value = '  -4.04%  '

$ws.Range("D33").Value = '168.84'
$ws.Range("E33").Value = '  +0.90%  '

$ws.Range("E34").Value = '  +4.32%  '

$ws.Range("D35").Value = '0.130'
$ws.Range("E35").Value = '  -1.04%  '

$ws.Range("D36").Value = '0.116'
$ws.Range("E36").Value = '  +3.30%  '

$ws.Range("D37").Value = '4.72'
$ws.Range("E37").Value = '  -0.67%  '

$ws.Range("D38").Value = '2.93'
$ws.Range("E38").Value = '  +3.31%  '

$ws.Range("D39").Value = '3.90'
$ws.Range("E39").Value = '  +1.93%  '

$ws.Range("D40").Value = '0.0353'
$ws.Range("E40").Value = '  -3.03%  '

$ws.Range("D41").Value = '1.74'
$ws.Range("E41").Value = '  +2.00%  '

$ws.Range("D42").Value = '99.58'
$ws.Range("E42").Value = '  -4.89%  '

$ws.Range("D43").Value = '0.229'
$ws.Range("E43").Value = '  -3.24%  '

$ws.Range("D44").Value = '69.14'
$ws.Range("E44").Value = '  -3.73%  '

$ws.Range("D45").Value = '12.84'
$ws.Range("E45").Value = '  -3.91%  '

$ws.Range("E46").Value = '  +0.42%  '

$ws.Range("D47").Value = '1.785.83'
$ws.Range("E47").Value = '  +7.13%  '

$ws.Range("D48").Value = '81.84'
$ws.Range("E48").Value = '  +8.23%  '

$ws.Range("D49").Value = '5.52'
$ws.Range("E49").Value = '  +3.57%  '

$ws.Range("D50").Value = '111.44'
$ws.Range("E50").Value = '  -2.45%  '

$ws.Range("D51").Value = '9.15'
$ws.Range("E51").Value = '  +2.23%  '
